$d = $word.ActiveDocument

# The first paragraph holds the "**ID__...__ID**" placeholder run.
$para = $d.Paragraphs.First

# Add a paragraph border (top/left/bottom/right) with 5pt distance from text,
# i.e. <w:pBdr><w:top w:space="5"/><w:left w:space="5"/><w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$para.Borders.DistanceFromTop = 5
$para.Borders.DistanceFromLeft = 5
$para.Borders.DistanceFromBottom = 5
$para.Borders.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$para.Format.LeftIndent = 11.25

# Swap the placeholder id text and drop the now-redundant trailing-space run
# by replacing "<old id><space>" with "<new id>" across the paragraph.
$para.Range.Find.Execute("**ID__AFFARS_pgi_5339_topic_2__ID** ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5339__ID**", 2)
